{"js": "// Apply the per-cell text replacements described by the diff: the header\n// date line, plus the 25 \"NNN\u00d7N=\" multiplication prompts scattered across\n// the single worksheet-style table (5 data rows x 5 columns).\nconst replacements = [\n  [\"2024-01-12 Friday\", \"2024-01-13 Saturday\"],\n  [\"626\u00d79=\", \"550\u00d79=\"],\n  [\"485\u00d76=\", \"612\u00d77=\"],\n  [\"304\u00d75=\", \"194\u00d75=\"],\n  [\"656\u00d78=\", \"232\u00d76=\"],\n  [\"266\u00d73=\", \"465\u00d75=\"],\n  [\"782\u00d75=\", \"533\u00d77=\"],\n  [\"177\u00d76=\", \"662\u00d75=\"],\n  [\"376\u00d77=\", \"737\u00d77=\"],\n  [\"715\u00d72=\", \"107\u00d75=\"],\n  [\"613\u00d77=\", \"359\u00d75=\"],\n  [\"750\u00d76=\", \"475\u00d77=\"],\n  [\"542\u00d73=\", \"554\u00d72=\"],\n  [\"126\u00d72=\", \"265\u00d77=\"],\n  [\"537\u00d75=\", \"738\u00d74=\"],\n  [\"559\u00d76=\", \"636\u00d79=\"],\n  [\"832\u00d76=\", \"393\u00d76=\"],\n  [\"500\u00d74=\", \"700\u00d73=\"],\n  [\"147\u00d77=\", \"241\u00d76=\"],\n  [\"952\u00d78=\", \"124\u00d78=\"],\n  [\"457\u00d73=\", \"321\u00d74=\"],\n  [\"904\u00d79=\", \"418\u00d72=\"],\n  [\"114\u00d74=\", \"934\u00d75=\"],\n  [\"747\u00d73=\", \"416\u00d73=\"],\n  [\"994\u00d75=\", \"365\u00d76=\"],\n  [\"172\u00d74=\", \"829\u00d77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the per-cell text replacements described by the diff: the header\n# date line, plus the 25 \"NNN\u00d7N=\" multiplication prompts scattered across\n# the single worksheet-style table (5 data rows x 5 columns).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-01-12 Friday\", \"2024-01-13 Saturday\"),\n    @(\"626\u00d79=\", \"550\u00d79=\"),\n    @(\"485\u00d76=\", \"612\u00d77=\"),\n    @(\"304\u00d75=\", \"194\u00d75=\"),\n    @(\"656\u00d78=\", \"232\u00d76=\"),\n    @(\"266\u00d73=\", \"465\u00d75=\"),\n    @(\"782\u00d75=\", \"533\u00d77=\"),\n    @(\"177\u00d76=\", \"662\u00d75=\"),\n    @(\"376\u00d77=\", \"737\u00d77=\"),\n    @(\"715\u00d72=\", \"107\u00d75=\"),\n    @(\"613\u00d77=\", \"359\u00d75=\"),\n    @(\"750\u00d76=\", \"475\u00d77=\"),\n    @(\"542\u00d73=\", \"554\u00d72=\"),\n    @(\"126\u00d72=\", \"265\u00d77=\"),\n    @(\"537\u00d75=\", \"738\u00d74=\"),\n    @(\"559\u00d76=\", \"636\u00d79=\"),\n    @(\"832\u00d76=\", \"393\u00d76=\"),\n    @(\"500\u00d74=\", \"700\u00d73=\"),\n    @(\"147\u00d77=\", \"241\u00d76=\"),\n    @(\"952\u00d78=\", \"124\u00d78=\"),\n    @(\"457\u00d73=\", \"321\u00d74=\"),\n    @(\"904\u00d79=\", \"418\u00d72=\"),\n    @(\"114\u00d74=\", \"934\u00d75=\"),\n    @(\"747\u00d73=\", \"416\u00d73=\"),\n    @(\"994\u00d75=\", \"365\u00d76=\"),\n    @(\"172\u00d74=\", \"829\u00d77=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    # 2 == wdReplaceAll\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n}\n"}
